# Update the cryptos price/volume table (columns D = Price, E = Volume(1h))
# with freshly scraped values, per the GitHub Actions scheduled refresh.
# NumberFormat is forced to "@" (Text) before each write so that Excel does
# not silently reinterpret values such as "1.000" or "251.99" as numbers
# (which would drop formatting like trailing zeros or thousands separators).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.352.03'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.57%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.942.57'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.46%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.99'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7259'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -8.08%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3354'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -4.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '28.86'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07400'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +5.85%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8225'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.76%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08145'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.940.58'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.512'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.45%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '95.46'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -4.87%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.92'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.82%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.363.98'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008359'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +5.88%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '254.34'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -6.74%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.880'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.196.61'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9994'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.998'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.944'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.81%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.34'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.414'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +4.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.42'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1320'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -11.42%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.31%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.345'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.480'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.270'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05324'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.308'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +6.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7639'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.750'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02000'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.39%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.846'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '81.52'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.89%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.624'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4583'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.050'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8455'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9995'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.16'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.66%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.857'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.535'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '37.17'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4232'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.521'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.04%  '
